$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above what is currently row 20 ("R12 ... CRCW0805-1.0M-E3 ..."),
# shifting the existing rows 20-26 (and the totals row 28) down by one.
$ws.Rows.Item(20).Insert()

# The newly inserted row inherited formatting from the row above (row 19); clear the
# cells we don't want populated so they serialize the same as a genuinely blank cell.
$ws.Cells.Item(20,1).Clear()
$ws.Cells.Item(20,2).Clear()
$ws.Cells.Item(20,4).Clear()
$ws.Cells.Item(20,5).Clear()
$ws.Cells.Item(20,8).Clear()
$ws.Cells.Item(20,10).Clear()

# Populate the new row: Distributor = "Mouser", Distributor PART # = "594-5073NW5R600J"
$ws.Cells.Item(20,6).Value = "Mouser"
$ws.Cells.Item(20,7).Value = "594-5073NW5R600J"

# Match the formatting used by the other "Distributor PART #" entries (e.g. G7), which
# use the Arial / dark-grey font style instead of the default.
$ws.Cells.Item(7,7).Copy()
$ws.Cells.Item(20,7).PasteSpecial(-4122)

# Restore the selected cell to what the author left it at.
$null = $ws.Range("H20").Select()
